$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Sending cluster" labels in column A
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("A5").Value = "Inflammatory-Mac"

# Row 2
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.724909
$ws.Range("H2").Value = 1.449818
$ws.Range("I2").Value = 0.5784209917375155
$ws.Range("J2").Value = 0.477722171991027
$ws.Range("O2").Value = 0.1797880856321904
$ws.Range("P2").Value = 0.2474388498825759
$ws.Range("Q2").Value = 0.07871955976433334
$ws.Range("R2").Value = 0.472317358586
$ws.Range("S2").Value = 0.1039932027939609
$ws.Range("T2").Value = 0.1182070248008658

# Row 3
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.724909
$ws.Range("H3").Value = 1.449818
$ws.Range("I3").Value = 0.5784209917375155
$ws.Range("J3").Value = 0.477722171991027
$ws.Range("M3").Value = 0.4954095
$ws.Range("N3").Value = 0.990819
$ws.Range("O3").Value = 0.8202119143678096
$ws.Range("P3").Value = 0.752561150117424
$ws.Range("Q3").Value = 0.3591268052355
$ws.Range("R3").Value = 1.436507220942
$ws.Range("S3").Value = 0.4744277889435546
$ws.Range("T3").Value = 0.3595151471901611

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.528346
$ws.Range("H4").Value = 1.585038
$ws.Range("I4").Value = 0.4215790082624845
$ws.Range("J4").Value = 0.522277828008973
$ws.Range("O4").Value = 0.1797880856321904
$ws.Range("P4").Value = 0.2474388498825759
$ws.Range("Q4").Value = 0.05737432494733333
$ws.Range("R4").Value = 0.5163689245259999
$ws.Range("S4").Value = 0.07579488283822945
$ws.Range("T4").Value = 0.1292318250817101

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.528346
$ws.Range("H5").Value = 1.585038
$ws.Range("I5").Value = 0.4215790082624845
$ws.Range("J5").Value = 0.522277828008973
$ws.Range("M5").Value = 0.4954095
$ws.Range("N5").Value = 0.990819
$ws.Range("O5").Value = 0.8202119143678096
$ws.Range("P5").Value = 0.752561150117424
$ws.Range("Q5").Value = 0.261747627687
$ws.Range("R5").Value = 1.570485766122
$ws.Range("S5").Value = 0.345784125424255
$ws.Range("T5").Value = 0.3930460029272629
